$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Recomputed strikeout ("K") counts per game row, written into column G
# (row number => new K value), replacing the previous values that were
# derived from the old "Strike#" stat.
$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 2
    6  = 1
    7  = 0
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 1
    33 = 0
    34 = 2
    35 = 1
    36 = 2
    37 = 1
    38 = 0
    39 = 1
    40 = 3
    41 = 0
    42 = 0
    43 = 1
    44 = 1
    45 = 1
    46 = 1
    47 = 1
    48 = 3
    49 = 1
    50 = 2
    51 = 2
    52 = 1
    53 = 1
    54 = 2
    55 = 1
    56 = 2
    57 = 1
    58 = 1
    59 = 0
    60 = 2
    61 = 1
    62 = 1
    64 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
